# "Before" workbook: calculate watt usage fix
# - Remove the stray/placeholder "fvjhtjhtjht" sheet.
# - Add the missing Drive + PSU components to the "jan" sheet.
# - Add the missing GPU component to the "tom" sheet.
# - Clear the leftover cell style on jan!B4.

$wb = $excel.ActiveWorkbook

# Remove the extraneous sheet that was left in the workbook.
$extra = $wb.Worksheets.Item("fvjhtjhtjht")
$extra.Delete()

# "tom" sheet: add the missing GPU row.
$tom = $wb.Worksheets.Item("tom")
$tom.Range("A4").Value = "GPU"
$tom.Range("B4").Value = "GTX Titan X"

# "jan" sheet: add the missing Drive and PSU rows.
$jan = $wb.Worksheets.Item("jan")
$jan.Range("A6").Value = "Drive"
$jan.Range("B6").Value = "950 EVO "
$jan.Range("A7").Value = "PSU"
$jan.Range("B7").Value = "G550M"

# Clear the stray style that was applied to B4 on the "jan" sheet.
$jan.Range("B4").ClearFormats()
